$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" '58.478.25'
Set-TextValue "E2" '  +1.72%  '
Set-TextValue "D3" '3.100.10'
Set-TextValue "E3" '  +0.34%  '
Set-TextValue "E4" '  -0.02%  '
Set-TextValue "D5" '527.77'
Set-TextValue "E5" '  +2.10%  '
Set-TextValue "D6" '142.92'
Set-TextValue "E6" '  +0.82%  '
Set-TextValue "E7" '  +0.01%  '
Set-TextValue "D8" '0.442'
Set-TextValue "E8" '  +1.47%  '
Set-TextValue "D9" '7.34'
Set-TextValue "E9" '  +0.80%  '
Set-TextValue "E10" '  +0.79%  '
Set-TextValue "E11" '  +2.38%  '
Set-TextValue "D12" '3.633.57'
Set-TextValue "E12" '  +0.50%  '
Set-TextValue "E13" '  +0.99%  '
Set-TextValue "D14" '26.90'
Set-TextValue "E14" '  +4.81%  '
Set-TextValue "D15" '0.0000167'
Set-TextValue "E15" '  +1.98%  '
Set-TextValue "D16" '58.530.09'
Set-TextValue "E16" '  +1.64%  '
Set-TextValue "D17" '3.096.89'
Set-TextValue "E17" '  +0.29%  '
Set-TextValue "E18" '  -0.55%  '
Set-TextValue "D19" '12.93'
Set-TextValue "E19" '  -2.34%  '
Set-TextValue "D20" '8.09'
Set-TextValue "E20" '  -0.81%  '
Set-TextValue "D21" '341.18'
Set-TextValue "E21" '  +1.89%  '
Set-TextValue "E22" '  -0.35%  '
Set-TextValue "E23" '  +0.71%  '
Set-TextValue "D24" '66.01'
Set-TextValue "E24" '  +0.14%  '
Set-TextValue "E25" '  +0.31%  '
Set-TextValue "D26" '1.00'
Set-TextValue "E26" '  +0.02%  '
Set-TextValue "D27" '0.0₃0919'
Set-TextValue "E27" '  +0.60%  '
Set-TextValue "E28" '  +3.78%  '
Set-TextValue "D29" '7.27'
Set-TextValue "E29" '  +1.65%  '
Set-TextValue "E30" '  +3.03%  '
Set-TextValue "E31" '  +3.72%  '
Set-TextValue "E32" '  +0.43%  '
Set-TextValue "D33" '154.19'
Set-TextValue "E33" '  -0.45%  '
Set-TextValue "D34" '4.67'
Set-TextValue "E34" '  +2.49%  '
Set-TextValue "E35" '  +3.10%  '
Set-TextValue "D36" '26.97'
Set-TextValue "E36" '  -3.59%  '
Set-TextValue "E37" '  +3.61%  '
Set-TextValue "D38" '0.0678'
Set-TextValue "E38" '  +0.38%  '
Set-TextValue "D39" '3.141.29'
Set-TextValue "E39" '  +0.37%  '
Set-TextValue "B40" 'Mantle'
Set-TextValue "C40" 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue "D40" '0.679'
Set-TextValue "E40" '  +1.04%  '
Set-TextValue "B41" 'Filecoin'
Set-TextValue "C41" 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue "D41" '3.89'
Set-TextValue "E41" '  +0.31%  '
Set-TextValue "E42" '  +0.04%  '
Set-TextValue "E43" '  +7.64%  '
Set-TextValue "D44" '1.00'
Set-TextValue "E44" '  -0.04%  '
Set-TextValue "D45" '2.290.56'
Set-TextValue "E45" '  -0.06%  '
Set-TextValue "E46" '  +0.52%  '
Set-TextValue "D47" '20.86'
Set-TextValue "E47" '  +4.07%  '
Set-TextValue "D48" '0.965'
Set-TextValue "E48" '  +2.54%  '
Set-TextValue "E49" '  +1.64%  '
Set-TextValue "D50" '268.67'
Set-TextValue "E50" '  +6.51%  '
Set-TextValue "D51" '0.744'
Set-TextValue "E51" '  +8.05%  '
